$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 362
$ws.Range("I28").Value = 349
$ws.Range("J28").Value = 366.33334
$ws.Range("K28").Value = 349
$ws.Range("L28").Value = 366.33334
$ws.Range("M28").Value = 136
$ws.Range("N28").Value = -1336.33334
$ws.Range("H33").Value = 367.33334
$ws.Range("I33").Value = 367.33334
$ws.Range("K33").Value = 367.33334
$ws.Range("M33").Value = -138.33334
$ws.Range("H62").Value = 5951.1816
$ws.Range("J62").Value = 8461.5
$ws.Range("L62").Value = 8461.5
$ws.Range("N62").Value = -9709.5
$ws.Range("H65").Value = 5951.1816
$ws.Range("J65").Value = 8461.5
$ws.Range("L65").Value = 42307.5
$ws.Range("N65").Value = -48547.5
$ws.Range("H92").Value = 791.7143
$ws.Range("I92").Value = 258.6
$ws.Range("J92").Value = 2124.5
$ws.Range("K92").Value = 258.6
$ws.Range("L92").Value = 2124.5
$ws.Range("M92").Value = 989.4
$ws.Range("N92").Value = -4620.5
$ws.Range("H98").Value = 4189.4443
$ws.Range("I98").Value = 2261
$ws.Range("J98").Value = 6600
$ws.Range("K98").Value = 2261
$ws.Range("L98").Value = 6600
$ws.Range("M98").Value = -763
$ws.Range("N98").Value = -9596
$ws.Range("H101").Value = 3495
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
$ws.Range("H113").Value = 3511.25
$ws.Range("I113").Value = 3511.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3511.25
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -257.25
$ws.Range("H122").Value = 4189.4443
$ws.Range("I122").Value = 2261
$ws.Range("J122").Value = 6600
$ws.Range("K122").Value = 6783
$ws.Range("L122").Value = 19800
$ws.Range("M122").Value = -4333
$ws.Range("N122").Value = -24700
$ws.Range("H138").Value = 2552.3572
$ws.Range("J138").Value = 3924.75
$ws.Range("L138").Value = 11774.25
$ws.Range("N138").Value = -22054.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8182.1665
$ws.Range("I32").Value = 8182.1665
$ws.Range("K32").Value = 8182.1665
$ws.Range("M32").Value = -7895.1665
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -9727
$ws.Range("H45").Value = 4756
$ws.Range("I45").Value = 4756
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4756
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4379
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H63").Value = 1915.6666
$ws.Range("I63").Value = 1874.75
$ws.Range("J63").Value = 1997.5
$ws.Range("K63").Value = 1874.75
$ws.Range("L63").Value = 1997.5
$ws.Range("M63").Value = -1188.75
$ws.Range("N63").Value = -3369.5
$ws.Range("H66").Value = 1915.6666
$ws.Range("I66").Value = 1874.75
$ws.Range("J66").Value = 1997.5
$ws.Range("K66").Value = 9373.75
$ws.Range("L66").Value = 9987.5
$ws.Range("M66").Value = -5941.75
$ws.Range("N66").Value = -16851.5
$ws.Range("H97").Value = 667.5
$ws.Range("I97").Value = 428.125
$ws.Range("J97").Value = 1625
$ws.Range("K97").Value = 428.125
$ws.Range("L97").Value = 1625
$ws.Range("M97").Value = 67.875
$ws.Range("N97").Value = -2617
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 35666.668
$ws.Range("I29").Value = 28500
$ws.Range("J29").Value = 50000
$ws.Range("K29").ClearContents()
$ws.Range("L29").Value = 50000
$ws.Range("M29").Value = -28211
$ws.Range("N29").Value = -50578
$ws.Range("H94").Value = 901.1667
$ws.Range("I94").Value = 867.3333
$ws.Range("J94").Value = 935
$ws.Range("K94").Value = 867.3333
$ws.Range("L94").Value = 935
$ws.Range("M94").Value = -416.3333
$ws.Range("N94").Value = -1837
$ws.Range("H134").Value = 3428.4
$ws.Range("I134").Value = 2443.4167
$ws.Range("K134").Value = 7330.250100000001
$ws.Range("M134").Value = -4795.250100000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 130000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H115").Value = 94999
$ws.Range("J115").Value = 94999
$ws.Range("L115").Value = 94999
$ws.Range("N115").Value = -97349
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 388
$ws.Range("I17").Value = 225
$ws.Range("K17").Value = 675
$ws.Range("M17").Value = -506
$ws.Range("H39").Value = 4233.3335
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H129").Value = 950
$ws.Range("I129").Value = 950
$ws.Range("K129").Value = 2850
$ws.Range("M129").Value = 2150
$ws.Range("H131").Value = 2000
$ws.Range("J131").Value = 2000
$ws.Range("L131").Value = 6000
$ws.Range("N131").Value = -16080
$ws.Range("H139").Value = 1523
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9977
$ws.Range("I80").Value = 3287.8
$ws.Range("J80").Value = 26700
$ws.Range("K80").Value = 3287.8
$ws.Range("L80").Value = 26700
$ws.Range("M80").Value = -2289.8
$ws.Range("N80").Value = -28696
$ws.Range("H83").Value = 9977
$ws.Range("I83").Value = 3287.8
$ws.Range("J83").Value = 26700
$ws.Range("K83").Value = 16439
$ws.Range("L83").Value = 133500
$ws.Range("M83").Value = -11447
$ws.Range("N83").Value = -143484
$ws.Range("H102").Value = 1496.3334
$ws.Range("I102").Value = 1496.3334
$ws.Range("K102").Value = 1496.3334
$ws.Range("M102").Value = 125.6666
$ws.Range("H107").Value = 666.6667
$ws.Range("I107").Value = 590.25
$ws.Range("K107").Value = 590.25
$ws.Range("M107").Value = 1329.75

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1507.5385
$ws.Range("I22").Value = 1287.25
$ws.Range("J22").Value = 1860
$ws.Range("K22").Value = 1287.25
$ws.Range("L22").Value = 1860
$ws.Range("M22").Value = -992.25
$ws.Range("N22").Value = -2450
$ws.Range("H27").Value = 1507.5385
$ws.Range("I27").Value = 1287.25
$ws.Range("J27").Value = 1860
$ws.Range("K27").Value = 1287.25
$ws.Range("L27").Value = 1860
$ws.Range("M27").Value = -1180.25
$ws.Range("N27").Value = -2074
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").ClearContents()
$ws.Range("N38").Value = 0
$ws.Range("H46").Value = 999.3333
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 999.5
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 999.5
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -1375.5
$ws.Range("H68").Value = 2898.5
$ws.Range("I68").Value = 2802
$ws.Range("K68").Value = 2802
$ws.Range("M68").Value = -2053
$ws.Range("H71").Value = 2898.5
$ws.Range("I71").Value = 2802
$ws.Range("K71").Value = 14010
$ws.Range("M71").Value = -10266
$ws.Range("H136").Value = 5376
$ws.Range("I136").Value = 5501.3335
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 16504.0005
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -13954.0005
$ws.Range("N136").Value = -20100

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 9997.5
$ws.Range("J48").Value = 9997.5
$ws.Range("L48").Value = 9997.5
$ws.Range("N48").Value = -11135.5
$ws.Range("H61").Value = 1500.6666
$ws.Range("I61").Value = 1501
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1501
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1209
$ws.Range("N61").Value = -2084
$ws.Range("H69").Value = 48390
$ws.Range("J69").Value = 48390
$ws.Range("L69").Value = 48390
$ws.Range("N69").Value = -49888
$ws.Range("H72").Value = 48390
$ws.Range("J72").Value = 48390
$ws.Range("L72").Value = 145170
$ws.Range("N72").Value = -152658
$ws.Range("H107").Value = 400
$ws.Range("I107").Value = 400
$ws.Range("K107").Value = 1200
$ws.Range("M107").Value = 720
$ws.Range("H122").Value = 2633.3333
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
